$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 6.735136143199226
$ws.Range("D2").Value = 9.853326502059996
$ws.Range("E2").Value = 12.11410413106104
$ws.Range("F2").Value = 55.15619820398764
$ws.Range("G2").Value = 3.593823875372472
$ws.Range("I2").Value = 17.47672228846301
$ws.Range("M2").Value = 28.55940843364026
$ws.Range("N2").Value = 16.62578978118533
$ws.Range("B3").Value = 6.577139174698647
$ws.Range("D3").Value = 9.495519704873601
$ws.Range("E3").Value = 11.50768105689166
$ws.Range("F3").Value = 52.75479987883801
$ws.Range("G3").Value = 3.606669334312419
$ws.Range("I3").Value = 17.60447137574947
$ws.Range("M3").Value = 27.2928102456399
$ws.Range("N3").Value = 16.67677300261363
$ws.Range("B4").Value = 6.478769257173619
$ws.Range("D4").Value = 9.275604784538192
$ws.Range("E4").Value = 11.11854256644482
$ws.Range("F4").Value = 51.25961235452363
$ws.Range("G4").Value = 3.614880172519094
$ws.Range("I4").Value = 17.68652627781072
$ws.Range("M4").Value = 26.49777230193483
$ws.Range("N4").Value = 16.71026101229107
$ws.Range("B5").Value = 6.438395636854663
$ws.Range("D5").Value = 9.18607387443479
$ws.Range("E5").Value = 10.95583054637351
$ws.Range("F5").Value = 50.6460683884268
$ws.Range("G5").Value = 3.618308622659099
$ws.Range("I5").Value = 17.72088030593627
$ws.Range("M5").Value = 26.16990630958212
$ws.Range("N5").Value = 16.72446108416407
$ws.Range("B6").Value = 6.431675937607436
$ws.Range("D6").Value = 9.171216862872042
$ws.Range("E6").Value = 10.92856553836618
$ws.Range("F6").Value = 50.54396401280863
$ws.Range("G6").Value = 3.618882925933351
$ws.Range("I6").Value = 17.72664028076113
$ws.Range("M6").Value = 26.11524518040367
$ws.Range("N6").Value = 16.7268525437938
$ws.Range("B7").Value = 6.478225852629073
$ws.Range("D7").Value = 9.274396796252089
$ws.Range("E7").Value = 11.11636478766642
$ws.Range("F7").Value = 51.25135370721622
$ws.Range("G7").Value = 3.614926074485251
$ws.Range("I7").Value = 17.68698587075284
$ws.Range("M7").Value = 26.49336562457787
$ws.Range("N7").Value = 16.71045027302651
$ws.Range("B8").Value = 6.680972866910538
$ws.Range("D8").Value = 9.730082215000845
$ws.Range("E8").Value = 11.90854149275889
$ws.Range("F8").Value = 54.33304253429402
$ws.Range("G8").Value = 3.598186469094532
$ws.Range("I8").Value = 17.52002360390755
$ws.Range("M8").Value = 28.1265531251949
$ws.Range("N8").Value = 16.64291804967674
$ws.Range("B9").Value = 7.065511738474994
$ws.Range("D9").Value = 10.61626551228617
$ws.Range("E9").Value = 13.32582650945784
$ws.Range("F9").Value = 60.17500217303715
$ws.Range("G9").Value = 3.567877998032225
$ws.Range("I9").Value = 17.22099864204126
$ws.Range("M9").Value = 31.17367061725379
$ws.Range("N9").Value = 16.52762441166351
$ws.Range("B10").Value = 7.337259030961415
$ws.Range("D10").Value = 11.25610178676127
$ws.Range("E10").Value = 14.28128078406007
$ws.Range("F10").Value = 64.3029134886602
$ws.Range("G10").Value = 3.547072620772043
$ws.Range("I10").Value = 17.01817015192901
$ws.Range("M10").Value = 33.29757803383404
$ws.Range("N10").Value = 16.45309327540301
$ws.Range("B11").Value = 7.458011726469839
$ws.Range("D11").Value = 11.54357482548543
$ws.Range("E11").Value = 14.69697741251307
$ws.Range("F11").Value = 66.13848210952544
$ws.Range("G11").Value = 3.537908646372871
$ws.Range("I11").Value = 16.92946146060132
$ws.Range("M11").Value = 34.23570583817392
$ws.Range("N11").Value = 16.42133566478197
$ws.Range("B12").Value = 7.503285426557074
$ws.Range("D12").Value = 11.6518370514433
$ws.Range("E12").Value = 14.85165264538756
$ws.Range("F12").Value = 66.82705934091932
$ws.Range("G12").Value = 3.53448028867665
$ws.Range("I12").Value = 16.89637339110338
$ws.Range("M12").Value = 34.58671466592187
$ws.Range("N12").Value = 16.40961344819153
$ws.Range("B13").Value = 7.49355567061288
$ws.Range("D13").Value = 11.62854857802216
$ws.Range("E13").Value = 14.81846258690563
$ws.Range("F13").Value = 66.67905764249851
$ws.Range("G13").Value = 3.53521680846824
$ws.Range("I13").Value = 16.90347721100795
$ws.Range("M13").Value = 34.51131005237062
$ws.Range("N13").Value = 16.41212461216806
$ws.Range("B14").Value = 7.461745714754129
$ws.Range("D14").Value = 11.5524940197542
$ws.Range("E14").Value = 14.70975761021141
$ws.Range("F14").Value = 66.19526450496906
$ws.Range("G14").Value = 3.537625762846044
$ws.Range("I14").Value = 16.92672923018505
$ws.Range("M14").Value = 34.26466955410903
$ws.Range("N14").Value = 16.42036521381437
$ws.Range("B15").Value = 7.4422010935122
$ws.Range("D15").Value = 11.5058283944099
$ws.Range("E15").Value = 14.64281552262175
$ws.Range("F15").Value = 65.89806753878172
$ws.Range("G15").Value = 3.539106724103554
$ws.Range("I15").Value = 16.9410371612963
$ws.Range("M15").Value = 34.1130375563018
$ws.Range("N15").Value = 16.42545221870516
$ws.Range("B16").Value = 7.329306337525645
$ws.Range("D16").Value = 11.23723525120751
$ws.Range("E16").Value = 14.25373125157809
$ws.Range("F16").Value = 64.18206529375638
$ws.Range("G16").Value = 3.547677439044635
$ws.Range("I16").Value = 17.02403842971437
$ws.Range("M16").Value = 33.23568616359744
$ws.Range("N16").Value = 16.45521148662777
$ws.Range("B17").Value = 7.259284999942365
$ws.Range("D17").Value = 11.07148039565662
$ws.Range("E17").Value = 14.01017424411825
$ws.Range("F17").Value = 63.11819482938965
$ws.Range("G17").Value = 3.553011303348935
$ws.Range("I17").Value = 17.07586310484129
$ws.Range("M17").Value = 32.6901153013173
$ws.Range("N17").Value = 16.47401420232071
$ws.Range("B18").Value = 7.218742981528475
$ws.Range("D18").Value = 10.97580855055878
$ws.Range("E18").Value = 13.86830387868837
$ws.Range("F18").Value = 62.50232937052888
$ws.Range("G18").Value = 3.556107570696202
$ws.Range("I18").Value = 17.10600675966006
$ws.Range("M18").Value = 32.37368686341707
$ws.Range("N18").Value = 16.48503142212625
$ws.Range("B19").Value = 7.204971479393479
$ws.Range("D19").Value = 10.94336101609416
$ws.Range("E19").Value = 13.81996344868237
$ws.Range("F19").Value = 62.29314339223452
$ws.Range("G19").Value = 3.55716082690096
$ws.Range("I19").Value = 17.11627073404899
$ws.Range("M19").Value = 32.26610482034211
$ws.Range("N19").Value = 16.48879657536906
$ws.Range("B20").Value = 7.266766890134198
$ws.Range("D20").Value = 11.089160482255
$ws.Range("E20").Value = 14.03628596270044
$ws.Range("F20").Value = 63.23185847152678
$ws.Range("G20").Value = 3.552440577968235
$ws.Range("I20").Value = 17.07031161519268
$ws.Range("M20").Value = 32.74846612662574
$ws.Range("N20").Value = 16.47199170989357
$ws.Range("B21").Value = 7.471101663379328
$ws.Range("D21").Value = 11.57484988862825
$ws.Range("E21").Value = 14.74176134236517
$ws.Range("F21").Value = 66.33754610058227
$ws.Range("G21").Value = 3.536917070242279
$ws.Range("I21").Value = 16.91988593394156
$ws.Range("M21").Value = 34.33723035569775
$ws.Range("N21").Value = 16.41793654800123
$ws.Range("B22").Value = 7.601990927522116
$ws.Range("D22").Value = 11.8887641150229
$ws.Range("E22").Value = 15.18685803725799
$ws.Range("F22").Value = 68.32918712379639
$ws.Range("G22").Value = 3.527014757751439
$ws.Range("I22").Value = 16.82450785739352
$ws.Range("M22").Value = 35.35079925764579
$ws.Range("N22").Value = 16.38437610861151
$ws.Range("B23").Value = 7.532388368635984
$ws.Range("D23").Value = 11.72156705721774
$ws.Range("E23").Value = 14.9507659345565
$ws.Range("F23").Value = 67.2698204314552
$ws.Range("G23").Value = 3.532278024951452
$ws.Range("I23").Value = 16.87514713656668
$ws.Range("M23").Value = 34.81216360656551
$ws.Range("N23").Value = 16.40212791007716
$ws.Range("B24").Value = 7.263385214729968
$ws.Range("D24").Value = 11.08116848360901
$ws.Range("E24").Value = 14.0244866039149
$ws.Range("F24").Value = 63.18048428689018
$ws.Range("G24").Value = 3.552698510120899
$ws.Range("I24").Value = 17.07282035614775
$ws.Range("M24").Value = 32.72209433805004
$ws.Range("N24").Value = 16.47290543268505
$ws.Range("B25").Value = 6.963183941754723
$ws.Range("D25").Value = 10.37806542212703
$ws.Range("E25").Value = 12.95733476459189
$ws.Range("F25").Value = 58.62083300311858
$ws.Range("G25").Value = 3.57581503438299
$ws.Range("I25").Value = 17.2988985685924
$ws.Range("M25").Value = 30.36827004104536
$ws.Range("N25").Value = 16.55700742504304
